$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: remove the now-unused "Preferred Callback Date/Time/Period" values ---
# (Preferred Callback Day "Tomorrow" in I2 is kept as-is)
$ws.Cells.Item(2, 7).ClearContents()   # G2
$ws.Cells.Item(2, 8).ClearContents()   # H2
$ws.Cells.Item(2, 10).ClearContents()  # J2

# --- Row 3: new call-queue entry ---
$ws.Cells.Item(3, 1).Value = "Aarav Mehta"
$ws.Cells.Item(3, 2).Value = "917823844614"
$ws.Cells.Item(3, 3).Value = "24 MG Road, Bengaluru"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "28"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(3, 5).Value = "Male"
$ws.Cells.Item(3, 6).Value = "2025-06-25 19:11:27"
$ws.Cells.Item(3, 11).Value = "Pending Callback"
$ws.Cells.Item(3, 12).Value = "Low"

# --- Row 4: new call-queue entry ---
$ws.Cells.Item(4, 1).Value = "Aarav Mehta"
$ws.Cells.Item(4, 2).Value = "917823844614"
$ws.Cells.Item(4, 3).Value = "24 MG Road, Bengaluru"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "28"
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).Value = "Male"
$ws.Cells.Item(4, 6).Value = "2025-06-25 19:35:25"
$ws.Cells.Item(4, 11).Value = "Pending Callback"
$ws.Cells.Item(4, 12).Value = "Low"

# --- Row 5: new call-queue entry ---
$ws.Cells.Item(5, 1).Value = "Aarav Mehta"
$ws.Cells.Item(5, 2).Value = "917823844614"
$ws.Cells.Item(5, 3).Value = "24 MG Road, Bengaluru"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "28"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "Male"
$ws.Cells.Item(5, 6).Value = "2025-06-25 19:39:22"
$ws.Cells.Item(5, 11).Value = "Pending Callback"
$ws.Cells.Item(5, 12).Value = "Low"

# --- Row 6: new call-queue entry ---
$ws.Cells.Item(6, 1).Value = "Vanshika panjwani"
$ws.Cells.Item(6, 2).Value = "917823844614"
$ws.Cells.Item(6, 3).Value = "24 MG Road, Bengaluru"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "28"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "Male"
$ws.Cells.Item(6, 6).Value = "2025-06-25 20:34:22"
$ws.Cells.Item(6, 7).Value = ""
$ws.Cells.Item(6, 8).Value = ""
$ws.Cells.Item(6, 9).Value = ""
$ws.Cells.Item(6, 10).Value = ""
$ws.Cells.Item(6, 11).Value = "Pending Callback"
$ws.Cells.Item(6, 12).Value = "Low"
